$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.527.21"
$ws.Range("E2").Value = "  +3.08%  "
$ws.Range("D3").Value = "2.996.30"
$ws.Range("E3").Value = "  +2.46%  "
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").Value = "563.12"
$ws.Range("E5").Value = "  +2.63%  "
$ws.Range("D6").Value = "138.95"
$ws.Range("E6").Value = "  +6.85%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("E8").Value = "  +2.11%  "
$ws.Range("D9").Value = "2.982.59"
$ws.Range("E9").Value = "  +2.20%  "
$ws.Range("E10").Value = "  +5.08%  "
$ws.Range("D11").Value = "5.27"
$ws.Range("E11").Value = "  +11.45%  "
$ws.Range("E12").Value = "  +2.00%  "
$ws.Range("D13").Value = "0.0000230"
$ws.Range("E13").Value = "  +5.00%  "
$ws.Range("D14").Value = "33.81"
$ws.Range("E14").Value = "  +3.73%  "
$ws.Range("E15").Value = "  -0.30%  "
$ws.Range("D16").Value = "3.490.35"
$ws.Range("E16").Value = "  +2.57%  "
$ws.Range("D17").Value = "7.18"
$ws.Range("D18").Value = "2.991.76"
$ws.Range("E18").Value = "  +2.45%  "
$ws.Range("D19").Value = "59.506.46"
$ws.Range("E19").Value = "  +3.18%  "
$ws.Range("D20").Value = "435.21"
$ws.Range("E20").Value = "  +4.72%  "
$ws.Range("E21").Value = "  +2.26%  "
$ws.Range("E22").Value = "  +3.99%  "
$ws.Range("D23").Value = "13.45"
$ws.Range("E23").Value = "  +1.11%  "
$ws.Range("D24").Value = "7.07"
$ws.Range("E24").Value = "  +1.66%  "
$ws.Range("D25").Value = "80.18"
$ws.Range("E25").Value = "  +0.97%  "
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  -0.07%  "
$ws.Range("E27").Value = "  +11.20%  "
$ws.Range("E28").Value = "  +0.26%  "
$ws.Range("E29").Value = "  +3.50%  "
$ws.Range("D30").Value = "7.77"
$ws.Range("E30").Value = "  +5.79%  "
$ws.Range("D31").Value = "6.24"
$ws.Range("E31").Value = "  +5.20%  "
$ws.Range("E32").Value = "  +2.41%  "
$ws.Range("D33").Value = "0.106"
$ws.Range("E33").Value = "  +9.54%  "
$ws.Range("D34").Value = "0.0₃0779"
$ws.Range("E34").Value = "  +13.27%  "
$ws.Range("E35").Value = "  +6.08%  "
$ws.Range("D36").Value = "5.88"
$ws.Range("E36").Value = "  +3.81%  "
$ws.Range("D37").Value = "2.09"
$ws.Range("E37").Value = "  +1.95%  "
$ws.Range("D38").Value = "48.92"
$ws.Range("E38").Value = "  +1.58%  "
$ws.Range("D39").Value = "8.59"
$ws.Range("E39").Value = "  -1.70%  "
$ws.Range("E40").Value = "  +7.67%  "
$ws.Range("D41").Value = "400.75"
$ws.Range("E41").Value = "  +7.89%  "
$ws.Range("E42").Value = "  +3.06%  "
$ws.Range("D43").Value = "2.759.23"
$ws.Range("E43").Value = "  +2.13%  "
$ws.Range("D44").Value = "0.106"
$ws.Range("E44").Value = "  -1.32%  "
$ws.Range("D45").Value = "0.251"
$ws.Range("E45").Value = "  +7.04%  "
$ws.Range("D47").Value = "122.89"
$ws.Range("E47").Value = "  -0.68%  "
$ws.Range("D48").Value = "34.45"
$ws.Range("E48").Value = "  +18.90%  "
$ws.Range("E49").Value = "  +1.93%  "
$ws.Range("E50").Value = "  +3.89%  "
$ws.Range("D51").Value = "23.54"
$ws.Range("E51").Value = "  +3.43%  "
